# Regenerate orders with updated distance/size codes.
# The experiment's stimulus naming scheme encodes a Distance code (D..)
# and a Size code (S..) inside Condition / Filename_Left / Filename_Right
# strings as well as the standalone Distance/Size lookup columns.
# This run renumbers:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# everywhere those tokens occur, leaving all other content (S25, S20,
# row/column layout, formatting, etc.) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Cells.Replace("D64", "D69")
[void]$ws.Cells.Replace("D80", "D86")
[void]$ws.Cells.Replace("D51", "D55")
[void]$ws.Cells.Replace("S30", "S31")
